# Added localized string for serviceworker modal
#
# Appends one new Key/String EN/String DE row to the "Tabelle2" table on
# the KeyValuePairs sheet:
#   Key        : service_worker-update_available
#   String EN  : A new update is available. Click the button below to
#                refresh the app and get the latest and greatest stuff!
#   String DE  : "\nEin neues Update ist verfügbar. Klicken Sie auf die
#                Schaltfläche unten, um die App zu aktualisieren und die
#                neuesten und besten Inhalte zu erhalten!"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Grow the table by one row (A1:C46 -> A1:C47). This also keeps the
# table's autoFilter ref and the sheet's <dimension> in sync.
$lo = $ws.ListObjects.Item(1)
$newRow = $lo.ListRows.Add()

$lastRow = $lo.Range.Rows.Count
$rowIdx = $lo.Range.Row + $lastRow - 1

# Match the formatting already used for ordinary (non-wrapped) vs.
# wrapped/multi-line table rows elsewhere in the sheet:
#  - column A keeps the plain fill-only style used throughout the table
#  - columns B/C use the wrap-text style used by other long entries
$ws.Range("A46").Copy() | Out-Null
$ws.Cells.Item($rowIdx, 1).PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("B14:C14").Copy() | Out-Null
$ws.Range($ws.Cells.Item($rowIdx, 2), $ws.Cells.Item($rowIdx, 3)).PasteSpecial(-4122) | Out-Null

$ws.Cells.Item($rowIdx, 1).Value = "service_worker-update_available"
$ws.Cells.Item($rowIdx, 2).Value = "A new update is available. Click the button below to refresh the app and get the latest and greatest stuff!"
$ws.Cells.Item($rowIdx, 3).Value = "`nEin neues Update ist verfügbar. Klicken Sie auf die Schaltfläche unten, um die App zu aktualisieren und die neuesten und besten Inhalte zu erhalten!"

# The longer bilingual strings wrap onto several lines, so the row needs
# to grow accordingly (same pattern as the existing wrapped rows).
$ws.Rows.Item($rowIdx).RowHeight = 57

# Reflect the scrolled-down view / new selection on the added row.
$ws.Cells.Item($rowIdx, 1).Select() | Out-Null
